$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.939.20'
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").Value = '1.638.41'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '215.44'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("D9").Value = '0.0640'
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = '19.64'
$ws.Range("E10").Value = '  -1.69%  '
$ws.Range("D11").Value = '0.0795'
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").Value = '1.865.16'
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("D13").Value = '4.25'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").Value = '1.632.45'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("E15").Value = '  -1.09%  '
$ws.Range("D16").Value = '0.0₃0765'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '62.97'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = '25.947.04'
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("D24").Value = '144.07'
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  +2.68%  '
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("D29").Value = '15.55'
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").Value = '3.31'
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("E34").Value = '  -3.62%  '
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("E36").Value = '  -1.40%  '
$ws.Range("D37").Value = '1.140.19'
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("D38").Value = '0.545'
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("D39").Value = '2.46'
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("D43").Value = '99.40'
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("D44").Value = '0.796'
$ws.Range("D45").Value = '1.774.97'
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("D47").Value = '56.65'
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("D48").Value = '0.0532'
$ws.Range("E48").Value = '  +2.88%  '
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("D50").Value = '7.66'
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("E51").Value = '  -0.83%  '
